$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Admittans (p.u.)") between the existing
# "Impedans (p.u.)" and "Kapasitans (nF)" columns, shifting Kapasitans to E.
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "Admittans (p.u.)"

# Row 2 (line 1-2)
$ws.Range("B2").Value = "(0.9159999999999999+9.847j)"
$ws.Range("D2").Value = "(0.8429235398105763-9.061428052963697j)"

# Row 3 (line 2-3)
$ws.Range("D3").Value = "(1.0267526096628827-11.03759055387599j)"

# Row 4 (line 3-4)
$ws.Range("D4").Value = "(0.4038273862272426-4.341144401942858j)"

# Row 5 (line 4-5)
$ws.Range("D5").Value = "(2.0322449532583664-22.35469448584203j)"

# Row 6 (line 5-6)
$ws.Range("B6").Value = "(1.4240000000000002+22.784000000000002j)"
$ws.Range("D6").Value = "(0.24592314082105537-3.934770253136886j)"

# Row 7 (line 6-7)
$ws.Range("D7").Value = "(0.3218699931334401-5.149919890135042j)"

# Row 8 (line 7-8)
$ws.Range("B8").Value = "(1.976+21.241999999999997j)"
$ws.Range("D8").Value = "(0.3907479567138098-4.200540534673455j)"

# Row 9 (line 1-8)
$ws.Range("D9").Value = "(0.182447533664104-1.9613109868891176j)"

# Row 10 (line 1-6)
$ws.Range("D10").Value = "(0.16540659007422623-1.7781208432979319j)"

# Re-apply bold to force the header font style to be re-serialized
# (adds family id to the bold font, matching the updated export).
$ws.Range("A1").Font.Bold = $true
$ws.Range("B1").Font.Bold = $true
$ws.Range("C1").Font.Bold = $true
$ws.Range("D1").Font.Bold = $true
$ws.Range("E1").Font.Bold = $true

# Column widths for the widened / new columns.
$ws.Columns.Item(2).ColumnWidth = 42.7505
$ws.Columns.Item(3).ColumnWidth = 46.417166666666674
$ws.Columns.Item(4).ColumnWidth = 42.917166666666674
$ws.Columns.Item(5).ColumnWidth = 14.583833333333335

# Selection as left by the author.
$ws.Range("C13").Select()
